# Reproduce the "fill range of columns, stop at existing data" feature result.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: fill J1:V1 with "OK"
$ws.Range("J1:V1").Value = "OK"

# Row 3: fill J3:N3 with "OK"; O3 already contains "p" (pre-existing data), so the fill stops there
$ws.Range("J3:N3").Value = "OK"
$ws.Range("O3").Value = "p"

# Row 5: B5 gets the same id value as A5
$ws.Range("B5").Value = 23760055

# Row 6: A6 now holds the numeric id instead of the text "valor"
$ws.Range("A6").Value = 23760055

# Rows 15-19: fill J:M with "OK"
$ws.Range("J15:M19").Value = "OK"

# Row 20: fill J20:M20 with "OK" (A20 already "OK")
$ws.Range("J20:M20").Value = "OK"

# Remove the now-obsolete OK markers in rows 26-40 (process stopped earlier this time)
$ws.Range("A26:A40").ClearContents()

# Leave the final selection where the feature's last operation left it
$ws.Range("O11").Select()
